# Scheduled-runner refresh of Yojimbo leve-profit cached prices (currentAveragePrice*
# / LevePrice* / LeveProfit* columns H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 8
$ws.Range("H8").Value = 75.44444
$ws.Range("I8").Value = 75.44444
$ws.Range("K8").Value = 226.33332
$ws.Range("M8").Value = -87.33332000000001

# row 112
$ws.Range("H112").Value = 1062.6097
$ws.Range("I112").Value = 718.7778
$ws.Range("J112").Value = 1105
$ws.Range("K112").Value = 2156.3334
$ws.Range("L112").Value = 3315
$ws.Range("M112").Value = -1048.3334
$ws.Range("N112").Value = -5531

# row 113
$ws.Range("H113").Value = 2237.2222
$ws.Range("I113").Value = 2104.375
$ws.Range("J113").Value = 3300
$ws.Range("K113").Value = 2104.375
$ws.Range("L113").Value = 3300
$ws.Range("M113").Value = 1149.625
$ws.Range("N113").Value = -9808

# row 137
$ws.Range("H137").Value = 3062.919
$ws.Range("I137").Value = 3054.1667
$ws.Range("J137").Value = 3100.4285
$ws.Range("K137").Value = 9162.500100000001
$ws.Range("L137").Value = 9301.2855
$ws.Range("M137").Value = -6612.500100000001
$ws.Range("N137").Value = -14401.2855

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 6096.55
$ws.Range("I32").Value = 5071.8335
$ws.Range("K32").Value = 5071.8335
$ws.Range("M32").Value = -4784.8335

$ws = $wb.Worksheets.Item("BSM")
# row 129
$ws.Range("H129").Value = 29528.908
$ws.Range("J129").Value = 29528.908
$ws.Range("L129").Value = 29528.908
$ws.Range("N129").Value = -39528.908

$ws = $wb.Worksheets.Item("CRP")
# row 14
$ws.Range("H14").Value = 60000
$ws.Range("J14").Value = 60000
$ws.Range("L14").Value = 60000
$ws.Range("N14").Value = -60340

# row 28
$ws.Range("H28").Value = 18643
$ws.Range("J28").Value = 18643
$ws.Range("L28").Value = 18643
$ws.Range("N28").Value = -19133

# row 53
$ws.Range("H53").Value = 16742
$ws.Range("J53").Value = 16742
$ws.Range("L53").Value = 16742
$ws.Range("N53").Value = -17956

# row 132
$ws.Range("H132").Value = 11505.667
$ws.Range("I132").Value = 10425.363
$ws.Range("J132").Value = 12694
$ws.Range("K132").Value = 31276.089
$ws.Range("L132").Value = 38082
$ws.Range("M132").Value = -28746.089
$ws.Range("N132").Value = -43142

# row 134
$ws.Range("H134").Value = 8035.375
$ws.Range("I134").Value = 2336
$ws.Range("J134").Value = 15363.143
$ws.Range("K134").Value = 7008
$ws.Range("L134").Value = 46089.429
$ws.Range("M134").Value = -4473
$ws.Range("N134").Value = -51159.429

$ws = $wb.Worksheets.Item("CUL")
# row 6
$ws.Range("H6").Value = 154.09091
$ws.Range("I6").Value = 105.75
$ws.Range("J6").Value = 283
$ws.Range("K6").Value = 317.25
$ws.Range("L6").Value = 849
$ws.Range("M6").Value = -204.25
$ws.Range("N6").Value = -1075

# row 32
$ws.Range("H32").Value = 4341.2
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 4341.2
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = 13023.6
$ws.Range("N32").Value = -13589.6

# row 39
$ws.Range("H39").Value = 2501.6667
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 2683.6365
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 8050.9095
$ws.Range("M39").Value = -1206
$ws.Range("N39").Value = -8638.9095

# row 46
$ws.Range("H46").Value = 2120.5
$ws.Range("I46").Value = 1058.5714
$ws.Range("J46").Value = 2692.3076
$ws.Range("K46").Value = 3175.7142
$ws.Range("L46").Value = 8076.9228
$ws.Range("M46").Value = -3084.7142
$ws.Range("N46").Value = -8258.9228

# row 58
$ws.Range("H58").Value = 6497.5
$ws.Range("I58").Value = 5000
$ws.Range("J58").Value = 7995
$ws.Range("K58").Value = 15000
$ws.Range("L58").Value = 23985
$ws.Range("M58").Value = -14872
$ws.Range("N58").Value = -24241

# row 76
$ws.Range("H76").Value = 10913.667
$ws.Range("I76").Value = 4500
$ws.Range("J76").Value = 12196.4
$ws.Range("K76").Value = 13500
$ws.Range("L76").Value = 36589.2
$ws.Range("M76").Value = -13117
$ws.Range("N76").Value = -37355.2

# row 79
$ws.Range("H79").Value = 10913.667
$ws.Range("I79").Value = 4500
$ws.Range("J79").Value = 12196.4
$ws.Range("K79").Value = 13500
$ws.Range("L79").Value = 36589.2
$ws.Range("M79").Value = -12174
$ws.Range("N79").Value = -39241.2

# row 98
$ws.Range("H98").Value = 625382.4
$ws.Range("I98").Value = 409.30768
$ws.Range("J98").Value = 3333599
$ws.Range("K98").Value = 1227.92304
$ws.Range("L98").Value = 10000797
$ws.Range("M98").Value = 270.0769599999999
$ws.Range("N98").Value = -10003793

# row 100
$ws.Range("H100").Value = 3169.5
$ws.Range("J100").Value = 3169.5
$ws.Range("L100").Value = 9508.5
$ws.Range("N100").Value = -11130.5

# row 106
$ws.Range("H106").Value = 3500
$ws.Range("J106").Value = 3500
$ws.Range("L106").Value = 10500
$ws.Range("N106").Value = -12392

# row 112
$ws.Range("H112").Value = 2777.3845
$ws.Range("I112").Value = 2209
$ws.Range("J112").Value = 3030
$ws.Range("K112").Value = 6627
$ws.Range("L112").Value = 9090
$ws.Range("M112").Value = -5519
$ws.Range("N112").Value = -11306

# row 118
$ws.Range("H118").Value = 2191.3547
$ws.Range("I118").Value = 672
$ws.Range("J118").Value = 2634.5
$ws.Range("K118").Value = 2016
$ws.Range("L118").Value = 7903.5
$ws.Range("M118").Value = -773
$ws.Range("N118").Value = -10389.5

# row 121
$ws.Range("H121").Value = 1378.2632
$ws.Range("I121").Value = 472.5
$ws.Range("J121").Value = 1619.8
$ws.Range("K121").Value = 1417.5
$ws.Range("L121").Value = 4859.4
$ws.Range("M121").Value = -107.5
$ws.Range("N121").Value = -7479.4

# row 122
$ws.Range("H122").Value = 1025.2667
$ws.Range("I122").Value = 432.85715
$ws.Range("J122").Value = 1543.625
$ws.Range("K122").Value = 3895.71435
$ws.Range("L122").Value = 13892.625
$ws.Range("M122").Value = -1445.71435
$ws.Range("N122").Value = -18792.625

# row 123
$ws.Range("H123").Value = 1278
$ws.Range("I123").Value = 516.6667
$ws.Range("J123").Value = 1398.2106
$ws.Range("K123").Value = 1550.0001
$ws.Range("L123").Value = 4194.6318
$ws.Range("M123").Value = 899.9999
$ws.Range("N123").Value = -9094.631799999999

# row 124
$ws.Range("H124").Value = 6400
$ws.Range("I124").Value = 2000
$ws.Range("J124").Value = 9333.333000000001
$ws.Range("K124").Value = 6000
$ws.Range("L124").Value = 27999.999
$ws.Range("M124").Value = -1090
$ws.Range("N124").Value = -37819.999

# row 125
$ws.Range("H125").Value = 4815
$ws.Range("I125").Value = 963.3333
$ws.Range("J125").Value = 8666.666999999999
$ws.Range("K125").Value = 2889.9999
$ws.Range("L125").Value = 26000.001
$ws.Range("M125").Value = 2030.0001
$ws.Range("N125").Value = -35840.001

# row 131
$ws.Range("H131").Value = 1115653.4
$ws.Range("I131").Value = 473.75
$ws.Range("J131").Value = 1261906.4
$ws.Range("K131").Value = 1421.25
$ws.Range("L131").Value = 3785719.2
$ws.Range("M131").Value = 3618.75
$ws.Range("N131").Value = -3795799.2

$ws = $wb.Worksheets.Item("GSM")
# row 9
$ws.Range("H9").Value = 6604.8
$ws.Range("I9").Value = 468
$ws.Range("J9").Value = 12741.6
$ws.Range("K9").Value = 468
$ws.Range("L9").Value = 12741.6
$ws.Range("M9").Value = -298
$ws.Range("N9").Value = -13081.6

# row 102
$ws.Range("H102").Value = 1367
$ws.Range("I102").Value = 1253.7333
$ws.Range("K102").Value = 1253.7333
$ws.Range("M102").Value = 368.2666999999999

# row 126
$ws.Range("H126").Value = 101146.2
$ws.Range("I126").Value = 200818.4
$ws.Range("J126").Value = 1474
$ws.Range("K126").Value = 602455.2
$ws.Range("L126").Value = 4422
$ws.Range("M126").Value = -599985.2
$ws.Range("N126").Value = -9362

# row 132
$ws.Range("H132").Value = 4220.048
$ws.Range("I132").Value = 4643.5483
$ws.Range("J132").Value = 3026.5454
$ws.Range("K132").Value = 13930.6449
$ws.Range("L132").Value = 9079.636200000001
$ws.Range("M132").Value = -11400.6449
$ws.Range("N132").Value = -14139.6362

$ws = $wb.Worksheets.Item("LTW")
# row 132
$ws.Range("H132").Value = 2419.5466
$ws.Range("I132").Value = 1992.0518
$ws.Range("J132").Value = 3878.0588
$ws.Range("K132").Value = 5976.1554
$ws.Range("L132").Value = 11634.1764
$ws.Range("M132").Value = -3446.1554
$ws.Range("N132").Value = -16694.1764

# row 133
$ws.Range("H133").Value = 28671.428
$ws.Range("J133").Value = 28671.428
$ws.Range("L133").Value = 28671.428
$ws.Range("N133").Value = -33731.428

# row 135
$ws.Range("H135").Value = 45799.832
$ws.Range("J135").Value = 45799.832
$ws.Range("L135").Value = 45799.832
$ws.Range("N135").Value = -55939.832

$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 2198.923
$ws.Range("I132").Value = 2349.9697
$ws.Range("J132").Value = 1368.1666
$ws.Range("K132").Value = 7049.909100000001
$ws.Range("L132").Value = 4104.4998
$ws.Range("M132").Value = -4519.909100000001
$ws.Range("N132").Value = -9164.4998
